$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Append 29 new rows (27-55) of "Capacity-Gen-Outage" download commands for
# operations, 2026, January 01-29 -- following the same layout/pattern as
# the existing RAMP-UNCERTAINTY rows (2-25).
# ---------------------------------------------------------------------------

$firstRow = 27
$lastRow = 55

# Columns A, B, C, D are constant for every new row.
$ws.Range("A${firstRow}:A${lastRow}").Value = "operations"
$ws.Range("B${firstRow}:B${lastRow}").Value = 2026
$ws.Range("D${firstRow}:D${lastRow}").Value = "Capacity-Gen-Outage"

# Column C keeps the text "01" (stored as text, quote-prefixed, like the
# existing month cells), and column E holds the yyyymmdd day stamp -- both
# vary only in that every row uses month "01" but a different day.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $day = $r - ($firstRow - 1)
    $ws.Range("C$r").Value = "'01"
    $dateStr = "202601{0:D2}" -f $day
    $ws.Range("E$r").Value = [long]$dateStr
}

# AI (source URL) and AJ (output path) formulas: row 27 is entered as a
# standalone formula, then rows 28-55 are entered as one fill so they share
# a single formula definition (mirrors how the original author typed the
# first new row, then filled the rest down).
$ws.Range("AI$firstRow").Formula = '=CONCAT("https://portal.spp.org/file-browser-api/download/capacity-of-generation-on-outage" & "?path=%2F" & B' + $firstRow + ' & "%2F" & C' + $firstRow + ' & "%2F" & D' + $firstRow + ' & "-" & E' + $firstRow + ' & ".csv")'
$ws.Range("AJ$firstRow").Formula = '=CONCAT("D:\pub_data_archive\" & LOWER(A' + $firstRow + ') &"\"& LOWER(D' + $firstRow + ') & "\" & B' + $firstRow + '&"\" & C' + $firstRow + ' & "\" & D' + $firstRow + ' &"-"&E' + $firstRow + '&".csv")'

$secondRow = $firstRow + 1
$ws.Range("AI${secondRow}:AI${lastRow}").Formula = '=CONCAT("https://portal.spp.org/file-browser-api/download/capacity-of-generation-on-outage" & "?path=%2F" & B' + $secondRow + ' & "%2F" & C' + $secondRow + ' & "%2F" & D' + $secondRow + ' & "-" & E' + $secondRow + ' & ".csv")'
$ws.Range("AJ${secondRow}:AJ${lastRow}").Formula = '=CONCAT("D:\pub_data_archive\" & LOWER(A' + $secondRow + ') &"\"& LOWER(D' + $secondRow + ') & "\" & B' + $secondRow + '&"\" & C' + $secondRow + ' & "\" & D' + $secondRow + ' &"-"&E' + $secondRow + '&".csv")'

# G (composed PowerShell command) formula for all new rows, filled as one
# shared definition across the whole new block.
$ws.Range("G${firstRow}:G${lastRow}").Formula = '=CONCAT("Start-Sleep -Seconds 3;" & " Invoke-WebRequest -Uri" & " ''" & AI' + $firstRow + ' & "''" & " -OutFile " & "''" & AJ' + $firstRow + ' & "'';")'

# Move the active selection to match where the author left off editing.
$ws.Range("E58").Select()
